$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value2 = 45309

# Step 2: update the computed prices in D14:D16
$ws.Range("D14").Value2 = 1266.597
$ws.Range("D15").Value2 = 1546.566
$ws.Range("D16").Value2 = 1817.002

# Re-apply the merged ranges in the new serialization order
$ws.Range("A1:D1").UnMerge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("A10:D10").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A31:D31").UnMerge()
$ws.Range("B13:C13").UnMerge()
$ws.Range("B14:C14").UnMerge()
$ws.Range("B15:C15").UnMerge()
$ws.Range("B16:C16").UnMerge()

$ws.Range("A10:D10").Merge()
$ws.Range("B15:C15").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("B13:C13").Merge()
$ws.Range("A1:D1").Merge()
$ws.Range("B16:C16").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("A31:D31").Merge()
$ws.Range("B14:C14").Merge()
